$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Website")

# New task rows - plain text (default style)
$ws.Range("D22").Value = "Payment gateway integration"
$ws.Range("D23").Value = "Card customization"

# New task row 27 - uses a red font style
$ws.Range("D27").Value = "If user input slash in keyword textbox"
$ws.Range("D27").Font.Color = 255
